# "Zips by Address File Group" - updating capstone date and other changes
#
# The SCF/MON "CBAT99" truckload group is being renamed to "CBAM99"
# (rows 47-83 and 156), and a single mis-bucketed zip (row 101, zip 11554)
# moves from the "CBAT10" truckload group to "CBAM8".
#
# Finally, the sheet's active selection is moved to I113 to match where the
# author was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename every "CBAT99" truckload bucket (column D) to "CBAM99".
# xlWhole (1) ensures this only matches full-cell values, so it does not
# clobber similarly-named buckets like "CBAT9".
$truckloadRange = $ws.Range("D1:D157")
$truckloadRange.Replace("CBAT99", "CBAM99", 1) | Out-Null

# Move zip 11554 (row 101) out of the "CBAT10" bucket into "CBAM8".
$ws.Range("D101").Value = "CBAM8"

# Reflect the author's final on-screen selection/scroll position.
$excel.Goto($ws.Range("I113"), $true)
